# Qatar Stars League workbook update (06-04-2024)
# This script:
#  1) Swaps the data (columns B:AC) between several pairs of adjacent rows
#     (the "id" column A and the "Date" column E stay where they are;
#     only the match-specific data that was mis-ordered moves).
#  2) Updates a handful of odds values for the two still-unplayed
#     fixtures in rows 112 and 113.
#  3) Appends two new upcoming fixtures as rows 114 and 115.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($sheet, [int]$row1, [int]$row2)
    $rng1 = $sheet.Range("B$row1" + ":AC$row1")
    $rng2 = $sheet.Range("B$row2" + ":AC$row2")
    $v1 = $rng1.Value()
    $v2 = $rng2.Value()
    $rng1.Value = $v2
    $rng2.Value = $v1
}

# --- 1) Row pair swaps -----------------------------------------------
Swap-RowData $ws 18 19
Swap-RowData $ws 27 28
Swap-RowData $ws 75 76
Swap-RowData $ws 87 88
Swap-RowData $ws 94 95
Swap-RowData $ws 96 97
Swap-RowData $ws 98 99
Swap-RowData $ws 108 109

# --- 2) Odds corrections for rows 112 & 113 ---------------------------
$ws.Cells.Item(112, 14).Value = 2.45    # N112 oddH
$ws.Cells.Item(112, 16).Value = 2.5     # P112 oddA
$ws.Cells.Item(112, 18).Value = 1.85    # R112 oddAHH
$ws.Cells.Item(112, 19).Value = 1.95    # S112 oddAHA

$ws.Cells.Item(113, 14).Value = 1.833   # N113 oddH
$ws.Cells.Item(113, 15).Value = 3.6     # O113 oddD
$ws.Cells.Item(113, 16).Value = 3.4     # P113 oddA
$ws.Cells.Item(113, 17).Value = -0.5    # Q113 Ah
$ws.Cells.Item(113, 18).Value = 1.9     # R113 oddAHH
$ws.Cells.Item(113, 19).Value = 1.9     # S113 oddAHA

# --- 3) Append new rows 114 and 115 -----------------------------------
# Copy the formatting (styles) from row 112, which has the same layout
# (no FTHG/FTAG/FTR yet, since these fixtures have not been played).
$ws.Range("A112:G112").Copy()
$ws.Range("A114:G114").PasteSpecial(-4122)
$ws.Range("A115:G115").PasteSpecial(-4122)

$ws.Range("K112:AA112").Copy()
$ws.Range("K114:AA114").PasteSpecial(-4122)
$ws.Range("K115:AA115").PasteSpecial(-4122)

# Row 114
$ws.Cells.Item(114, 1).Value = 112
$ws.Cells.Item(114, 2).Value = 7004659
$ws.Cells.Item(114, 3).Value = "Qatar Stars League"
$ws.Cells.Item(114, 4).Value = "Qatar Stars League"
$ws.Cells.Item(114, 5).Value = 45389.64583333334
$ws.Cells.Item(114, 6).Value = "AlRayyan SC"
$ws.Cells.Item(114, 7).Value = "AlMuaidar"
$ws.Cells.Item(114, 11).Value = 1.5
$ws.Cells.Item(114, 12).Value = 3.6
$ws.Cells.Item(114, 13).Value = 6
$ws.Cells.Item(114, 14).Value = 1.571
$ws.Cells.Item(114, 15).Value = 3.6
$ws.Cells.Item(114, 16).Value = 5.25
$ws.Cells.Item(114, 17).Value = -0.75
$ws.Cells.Item(114, 18).Value = 1.725
$ws.Cells.Item(114, 19).Value = 2.075
$ws.Cells.Item(114, 20).Value = 3.25
$ws.Cells.Item(114, 21).Value = 1.875
$ws.Cells.Item(114, 22).Value = 1.925
$ws.Cells.Item(114, 23).Value = 0
$ws.Cells.Item(114, 24).Value = 0
$ws.Cells.Item(114, 25).Value = 0
$ws.Cells.Item(114, 26).Value = 0
$ws.Cells.Item(114, 27).Value = 0

# Row 115
$ws.Cells.Item(115, 1).Value = 113
$ws.Cells.Item(115, 2).Value = 7004658
$ws.Cells.Item(115, 3).Value = "Qatar Stars League"
$ws.Cells.Item(115, 4).Value = "Qatar Stars League"
$ws.Cells.Item(115, 5).Value = 45389.64583333334
$ws.Cells.Item(115, 6).Value = "Al Gharafa"
$ws.Cells.Item(115, 7).Value = "Umm Salal"
$ws.Cells.Item(115, 11).Value = 1.571
$ws.Cells.Item(115, 12).Value = 3.25
$ws.Cells.Item(115, 13).Value = 6
$ws.Cells.Item(115, 14).Value = 1.727
$ws.Cells.Item(115, 15).Value = 3.25
$ws.Cells.Item(115, 16).Value = 4.5
$ws.Cells.Item(115, 17).Value = -0.75
$ws.Cells.Item(115, 18).Value = 1.85
$ws.Cells.Item(115, 19).Value = 1.95
$ws.Cells.Item(115, 20).Value = 3.25
$ws.Cells.Item(115, 21).Value = 1.85
$ws.Cells.Item(115, 22).Value = 1.95
$ws.Cells.Item(115, 23).Value = 0
$ws.Cells.Item(115, 24).Value = 0
$ws.Cells.Item(115, 25).Value = 0
$ws.Cells.Item(115, 26).Value = 0
$ws.Cells.Item(115, 27).Value = 0

Write-Output "done"
